# Add Ad and AdManager
# - B10: append "+0.5" to the hours worked ("2+1" -> "2+1+0.5")
# - C10: append note about the new AdsManager class
# - Move the active selection to C11 (next free row, Completed column)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "2+1+0.5"
$ws.Range("C10").Value = "Load images into Photo Wall. Setup the backend for ads. Add AdsManager to handle loading from Parse.om"

$ws.Range("C11").Select()
